# feat: add 2022-Q1 data
#
# 1. Clone the "2021-Q4" sheet (which already has the right layout/styles)
#    and drop it in right after "2021-Q4" / right before "总计", then
#    rename it to "2022-Q1" and fill in the new quarter's numbers.
# 2. Insert a new top data row into "总计" for 2022-Q1 and push the
#    existing history rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force the cell to stay text (the source data stores numeric-looking
    # strings like "2.46" as text, not numbers) without leaving a stray
    # number-format style behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "2022-Q1"

Set-TextValue $newSheet.Cells.Item(2, 4) "2.46"
Set-TextValue $newSheet.Cells.Item(2, 5) "81.85"
Set-TextValue $newSheet.Cells.Item(2, 6) "2.58"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.0635"
$newSheet.Cells.Item(2, 8).Value = 9

Set-TextValue $newSheet.Cells.Item(3, 4) "2.46"
Set-TextValue $newSheet.Cells.Item(3, 5) "81.85"
Set-TextValue $newSheet.Cells.Item(3, 6) "2.58"
Set-TextValue $newSheet.Cells.Item(3, 7) "0.0635"
$newSheet.Cells.Item(3, 8).Value = 9

# Restore the original active-sheet selection (cloning a sheet makes it
# the active one, which would otherwise shift tabSelected/activeTab).
$wb.Worksheets.Item(1).Activate()

# ---------------------------------------------------------------------
# 2. Prepend a 2022-Q1 row to the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the brand-new last row (7) the same column-A style as the row
# above it before shifting data down into it.
$total.Cells.Item(6, 1).Copy()
$total.Cells.Item(7, 1).PasteSpecial(-4122)

for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($r, 4).Value2
}

# Column A is just a fresh running row index (0, 1, 2, ...), not shifted data.
for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.13
